$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TFM_FILL marker cell gets "turned off" (prefixed with ~) ---
# U24 drives a chain of formulas (Z24 -> AE24 -> AJ24 -> AP24 -> AU24) that
# will pick up the new text automatically on recalculation.
$ws.Range("U24").Value = "~TFM_FILL"

# --- Updated Reference cost values (2017 / 2020 vintages) pulled in from the model ---
# Column X = 2017 cost, Column AC = 2020 cost (several technologies).
$ws.Range("X33").Value  = 81278.984337530506
$ws.Range("AC33").Value = 81279.984337530506

$ws.Range("X34").Value  = 48762.532743362797
$ws.Range("AC34").Value = 19505.013097345101

$ws.Range("X35").Value  = 51982.699999999997
$ws.Range("AC35").Value = 20793.080000000002

$ws.Range("X36").Value  = 32295.527999999998
$ws.Range("AC36").Value = 21530.351999999999

$ws.Range("AC44").Value = 38195.713600000003

$ws.Range("X45").Value  = 56564.537982300899
$ws.Range("AC45").Value = 22625.8151929204

$ws.Range("X46").Value  = 56564.537982300899
$ws.Range("AC46").Value = 22625.8151929204

$ws.Range("X47").Value  = 56564.537982300899
$ws.Range("AC47").Value = 22625.8151929204

$ws.Range("X48").Value  = 85334.432300885004
$ws.Range("AC48").Value = 34133.772920354

$ws.Range("X49").Value  = 56564.537982300899
$ws.Range("AC49").Value = 22625.8151929204

$ws.Range("X50").Value  = 1637.60737116

$ws.Range("X51").Value  = 337.44

$ws.Range("X52").Value  = 19046.524215158599

$ws.Range("X53").Value  = 9523.2621384296508

$ws.Range("X54").Value  = 4526.8304959535499

$ws.Range("X55").Value  = 4526.8304959535499

$ws.Range("X56").Value  = 7272.9396030149201

$ws.Range("X57").Value  = 4526.8304959535499

$ws.Range("X58").Value  = 4526.8304959535499

$ws.Range("X59").Value  = 4526.8304959535499

$ws.Range("X60").Value  = 4526.8304959535499

$ws.Range("X61").Value  = 7272.9396030149201

$ws.Range("X62").Value  = 7272.9396030149201

$ws.Range("X63").Value  = 19046.524215158599

$ws.Range("X64").Value  = 9523.2621384296508

$ws.Range("X65").Value  = 4761.6310692148199

# --- Leave the cursor on U25, matching where editing finished ---
$ws.Range("U25").Select()
